$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows right after row 30 (before the existing row 31) ---
# This pushes the existing rows 31-39 down to rows 33-41.
$ws.Rows.Item(31).Insert()
$ws.Rows.Item(31).Insert()

# Fill new row 31 with fresh weekly data
$ws.Range("A31").Value = 1
$ws.Range("B31").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C31").Value = "Arica y Parinacota"
$ws.Range("D31").Value = 44452
$ws.Range("E31").Value = 15
$ws.Range("F31").Value = 100112021
$ws.Range("G31").Value = "Ají"
$ws.Range("H31").Value = "Cristal"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 120
$ws.Range("K31").Value = 37000
$ws.Range("L31").Value = 38000
$ws.Range("M31").Value = 37500
$ws.Range("N31").Value = "$/caja 15 kilos"
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 2500
$ws.Range("Q31").Value = 15
$ws.Range("R31").Value = "Hortaliza"

# Fill new row 32 with fresh weekly data
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C32").Value = "Arica y Parinacota"
$ws.Range("D32").Value = 44452
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = 100112021
$ws.Range("G32").Value = "Ají"
$ws.Range("H32").Value = "Inferno"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 120
$ws.Range("K32").Value = 38000
$ws.Range("L32").Value = 40000
$ws.Range("M32").Value = 39000
$ws.Range("N32").Value = "$/caja 15 kilos"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 2600
$ws.Range("Q32").Value = 15
$ws.Range("R32").Value = "Hortaliza"

# --- Insert one new row before the old row 37 (now at row 39) ---
# This pushes rows 39-41 (previously 37-39) down to rows 40-42.
$ws.Rows.Item(39).Insert()

# Fill new row 39 with fresh weekly data
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C39").Value = "Arica y Parinacota"
$ws.Range("D39").Value = 44449
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = 100112021
$ws.Range("G39").Value = "Ají"
$ws.Range("H39").Value = "Inferno"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 120
$ws.Range("K39").Value = 40000
$ws.Range("L39").Value = 43000
$ws.Range("M39").Value = 41500
$ws.Range("N39").Value = "$/caja 15 kilos"
$ws.Range("O39").Value = "Región de Arica y Parinacota"
$ws.Range("P39").Value = 2767
$ws.Range("Q39").Value = 15
$ws.Range("R39").Value = "Hortaliza"
